# Attendance report logic fix - update Lieu From/To and Lieu Sum values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Lieu From/To now real times, Lieu Sum recalculated to 8.5
$ws.Range("L2").Value = "08:00"
$ws.Range("M2").Value = "18:00"
$ws.Range("N2").Value = 8.5

# Row 3: same fix as row 2
$ws.Range("L3").Value = "08:00"
$ws.Range("M3").Value = "18:00"
$ws.Range("N3").Value = 8.5

# Row 5: Lieu Sum corrected from 1.67 to 2.33, stored as text like the
# surrounding "pls choose"/placeholder cells in this column elsewhere in
# the sheet. Force text storage (apostrophe alone gets silently dropped by
# this host unless the cell's number format is Text), then restore the
# default "Normal" style so no stray formatting is left behind.
$ws.Range("N5").NumberFormat = "@"
$ws.Range("N5").Value = "2.33"
$ws.Range("N5").Style = "Normal"

# Row 16: Lieu Sum text trimmed from "2:" to "2" (still text)
$ws.Range("N16").NumberFormat = "@"
$ws.Range("N16").Value = "2"
$ws.Range("N16").Style = "Normal"
